$wb = $excel.ActiveWorkbook

# ALC!row139
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 25000
$ws.Range("J139").Value = 25000
$ws.Range("L139").Value = 25000
$ws.Range("N139").Value = -35280

# ARM!row4
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

# ARM!row104
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H104").Value = 49225
$ws.Range("J104").Value = 49225
$ws.Range("L104").Value = 49225
$ws.Range("N104").Value = -56213

# ARM!row110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 873.7273
$ws.Range("I110").Value = 873.7273
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 873.7273
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1171.2727
$ws.Range("N110").ClearContents()

# ARM!row122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 5129707
$ws.Range("I122").Value = 5129707
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 15389121
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -15386671
$ws.Range("N122").ClearContents()

# BSM!row64
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 500
$ws.Range("J64").Value = 500
$ws.Range("L64").Value = 500
$ws.Range("N64").Value = -950

# BSM!row67
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 500
$ws.Range("J67").Value = 500
$ws.Range("L67").Value = 500
$ws.Range("N67").Value = -2060

# BSM!row86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2292.8125
$ws.Range("I86").Value = 2189.8333
$ws.Range("J86").Value = 2601.75
$ws.Range("K86").Value = 2189.8333
$ws.Range("L86").Value = 2601.75
$ws.Range("M86").Value = -1066.8333
$ws.Range("N86").Value = -4847.75

# BSM!row89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2292.8125
$ws.Range("I89").Value = 2189.8333
$ws.Range("J89").Value = 2601.75
$ws.Range("K89").Value = 10949.1665
$ws.Range("L89").Value = 13008.75
$ws.Range("M89").Value = -5333.166499999999
$ws.Range("N89").Value = -24240.75

# BSM!row107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1127.1904
$ws.Range("I107").Value = 1045.8948
$ws.Range("J107").Value = 1899.5
$ws.Range("K107").Value = 1045.8948
$ws.Range("L107").Value = 1899.5
$ws.Range("M107").Value = 874.1052
$ws.Range("N107").Value = -5739.5

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8812.5625
$ws.Range("I31").Value = 1527.1765
$ws.Range("J31").Value = 17069.334
$ws.Range("K31").Value = 1527.1765
$ws.Range("L31").Value = 17069.334
$ws.Range("M31").Value = -1232.1765
$ws.Range("N31").Value = -17659.334

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 8812.5625
$ws.Range("I34").Value = 1527.1765
$ws.Range("J34").Value = 17069.334
$ws.Range("K34").Value = 1527.1765
$ws.Range("L34").Value = 17069.334
$ws.Range("M34").Value = -1325.1765
$ws.Range("N34").Value = -17473.334

# CRP!row132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1846.5807
$ws.Range("I132").Value = 1408.9231
$ws.Range("J132").Value = 4122.4
$ws.Range("K132").Value = 4226.7693
$ws.Range("L132").Value = 12367.2
$ws.Range("M132").Value = -1696.7693
$ws.Range("N132").Value = -17427.2

# CUL!row55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 3907
$ws.Range("J55").Value = 3907
$ws.Range("L55").Value = 11721
$ws.Range("N55").Value = -12075

# CUL!row68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2122.2222
$ws.Range("I68").Value = 850
$ws.Range("K68").Value = 2550
$ws.Range("M68").Value = -1739

# CUL!row71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 2122.2222
$ws.Range("I71").Value = 850
$ws.Range("K71").Value = 7650
$ws.Range("M71").Value = -3594

# CUL!row82
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 4000
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()

# CUL!row85
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H85").Value = 4000
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()

# CUL!row88
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 3934.7827
$ws.Range("J88").Value = 3934.7827
$ws.Range("L88").Value = 11804.3481
$ws.Range("N88").Value = -12660.3481

# CUL!row91
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H91").Value = 3934.7827
$ws.Range("J91").Value = 3934.7827
$ws.Range("L91").Value = 11804.3481
$ws.Range("N91").Value = -14768.3481

# CUL!row103
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 2566.1
$ws.Range("I103").Value = 397
$ws.Range("J103").Value = 3495.7144
$ws.Range("K103").Value = 1191
$ws.Range("L103").Value = 10487.1432
$ws.Range("M103").Value = -312
$ws.Range("N103").Value = -12245.1432

# CUL!row105
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H105").Value = 7800
$ws.Range("J105").Value = 7800
$ws.Range("L105").Value = 23400
$ws.Range("N105").Value = -28642

# CUL!row122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 801
$ws.Range("J122").Value = 850
$ws.Range("L122").Value = 7650
$ws.Range("N122").Value = -12550

# CUL!row126
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value = 4885.4546
$ws.Range("J126").Value = 4885.4546
$ws.Range("L126").Value = 14656.3638
$ws.Range("N126").Value = -24536.3638

# CUL!row130
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 4525
$ws.Range("I130").Value = 1000
$ws.Range("J130").Value = 4845.4546
$ws.Range("K130").Value = 3000
$ws.Range("L130").Value = 14536.3638
$ws.Range("M130").Value = 2020
$ws.Range("N130").Value = -24576.3638

# GSM!row70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5308.074
$ws.Range("I70").Value = 5222.951
$ws.Range("J70").Value = 5576.5386
$ws.Range("K70").Value = 5222.951
$ws.Range("L70").Value = 5576.5386
$ws.Range("M70").Value = -4952.951
$ws.Range("N70").Value = -6116.5386

# GSM!row73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5308.074
$ws.Range("I73").Value = 5222.951
$ws.Range("J73").Value = 5576.5386
$ws.Range("K73").Value = 5222.951
$ws.Range("L73").Value = 5576.5386
$ws.Range("M73").Value = -4286.951
$ws.Range("N73").Value = -7448.5386

# GSM!row113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 62501256
$ws.Range("I113").Value = 111112190
$ws.Range("J113").Value = 1480.4286
$ws.Range("K113").Value = 111112190
$ws.Range("L113").Value = 1480.4286
$ws.Range("M113").Value = -111110020
$ws.Range("N113").Value = -5820.4286

# GSM!row119
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H119").Value = 31000
$ws.Range("J119").Value = 31000
$ws.Range("L119").Value = 31000
$ws.Range("N119").Value = -40676

# GSM!row132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2302.8027
$ws.Range("I132").Value = 1923.52
$ws.Range("K132").Value = 5770.559999999999
$ws.Range("M132").Value = -3240.559999999999

# GSM!row136
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 39000
$ws.Range("J136").Value = 39000
$ws.Range("L136").Value = 117000
$ws.Range("N136").Value = -122100

# GSM!row139
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# LTW!row140
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 83286
$ws.Range("J140").Value = 83286
$ws.Range("L140").Value = 83286
$ws.Range("N140").Value = -93646

# WVR!row46
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 36940.668
$ws.Range("J46").Value = 36940.668
$ws.Range("L46").Value = 36940.668
$ws.Range("N46").Value = -37402.668

# WVR!row122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2758
$ws.Range("I122").Value = 1280.3
$ws.Range("J122").Value = 6452.25
$ws.Range("K122").Value = 3840.9
$ws.Range("L122").Value = 19356.75
$ws.Range("M122").Value = -1390.9
$ws.Range("N122").Value = -24256.75

# WVR!row134
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H134").Value = 36940.668
$ws.Range("J134").Value = 36940.668
$ws.Range("L134").Value = 110822.004
$ws.Range("N134").Value = -115892.004

# WVR!row137
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
